# Add 2022-Q4 data:
#  - insert a new "2022-Q4" sheet (positioned right after "总计", before "2022-Q3")
#  - populate it with the Q4 fund holdings
#  - update the "总计" summary sheet with the new quarter's row

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q4" worksheet by duplicating "2022-Q3" (so it keeps
#    the same header row / column styling), placed before "2022-Q3".
# ---------------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Item("2022-Q3")
$wsQ3.Copy($wsQ3)
$wsQ4 = $wb.Worksheets.Item(2)
$wsQ4.Name = "2022-Q4"

# Build rows 3-6 from row 2's existing formatting (copy limited to A:H so we
# don't drag whole-row formatting across all 16384 columns).
$wsQ4.Range("A2:H2").Copy($wsQ4.Range("A3:H3"))
$wsQ4.Range("A2:H2").Copy($wsQ4.Range("A4:H4"))
$wsQ4.Range("A2:H2").Copy($wsQ4.Range("A5:H5"))
$wsQ4.Range("A2:H2").Copy($wsQ4.Range("A6:H6"))

# Force the text-valued columns (B,C,D,E,F,G) to stay text even though some
# look numeric ("2.95", "012920", ...).
foreach ($r in 2..6) {
    foreach ($c in 2..7) {
        $wsQ4.Cells.Item($r, $c).NumberFormat = "@"
    }
}

# Row 2: 012920
$wsQ4.Cells.Item(2,1).Value = 0
$wsQ4.Cells.Item(2,2).Value = "012920"
$wsQ4.Cells.Item(2,3).Value = "易方达全球成长精选混合（QDII）A（人民币份额）"
$wsQ4.Cells.Item(2,4).Value = "2.95"
$wsQ4.Cells.Item(2,5).Value = "82.00"
$wsQ4.Cells.Item(2,6).Value = "1.78"
$wsQ4.Cells.Item(2,7).Value = "0.0525"
$wsQ4.Cells.Item(2,8).Value = 10

# Row 3: 012921
$wsQ4.Cells.Item(3,1).Value = 1
$wsQ4.Cells.Item(3,2).Value = "012921"
$wsQ4.Cells.Item(3,3).Value = "易方达全球成长精选混合（QDII）A（美元现汇份额）"
$wsQ4.Cells.Item(3,4).Value = "2.95"
$wsQ4.Cells.Item(3,5).Value = "82.00"
$wsQ4.Cells.Item(3,6).Value = "1.78"
$wsQ4.Cells.Item(3,7).Value = "0.0525"
$wsQ4.Cells.Item(3,8).Value = 10

# Row 4: 012922
$wsQ4.Cells.Item(4,1).Value = 2
$wsQ4.Cells.Item(4,2).Value = "012922"
$wsQ4.Cells.Item(4,3).Value = "易方达全球成长精选混合（QDII）C（人民币份额）"
$wsQ4.Cells.Item(4,4).Value = "2.95"
$wsQ4.Cells.Item(4,5).Value = "82.00"
$wsQ4.Cells.Item(4,6).Value = "1.78"
$wsQ4.Cells.Item(4,7).Value = "0.0525"
$wsQ4.Cells.Item(4,8).Value = 10

# Row 5: 012923
$wsQ4.Cells.Item(5,1).Value = 3
$wsQ4.Cells.Item(5,2).Value = "012923"
$wsQ4.Cells.Item(5,3).Value = "易方达全球成长精选混合（QDII）C（美元现汇份额）"
$wsQ4.Cells.Item(5,4).Value = "2.95"
$wsQ4.Cells.Item(5,5).Value = "82.00"
$wsQ4.Cells.Item(5,6).Value = "1.78"
$wsQ4.Cells.Item(5,7).Value = "0.0525"
$wsQ4.Cells.Item(5,8).Value = 10

# Row 6: 005698 (updated figures for this quarter)
$wsQ4.Cells.Item(6,1).Value = 4
$wsQ4.Cells.Item(6,2).Value = "005698"
$wsQ4.Cells.Item(6,3).Value = "华夏全球科技先锋混合（QDII）"
$wsQ4.Cells.Item(6,4).Value = "0.60"
$wsQ4.Cells.Item(6,5).Value = "83.35"
$wsQ4.Cells.Item(6,6).Value = "7.13"
$wsQ4.Cells.Item(6,7).Value = "0.0428"
$wsQ4.Cells.Item(6,8).Value = 3

# ---------------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet with the new quarter row.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# Make room for the extra row (5th table row), copying formatting down.
$wsTotal.Range("A4:D4").Copy($wsTotal.Range("A5:D5"))

$wsTotal.Cells.Item(2,1).Value = 0
$wsTotal.Cells.Item(2,2).Value = "2022-Q4"
$wsTotal.Cells.Item(2,3).Value = 5
$wsTotal.Cells.Item(2,4).Value = 0.25

$wsTotal.Cells.Item(3,1).Value = 1
$wsTotal.Cells.Item(3,2).Value = "2022-Q3"
$wsTotal.Cells.Item(3,3).Value = 1
$wsTotal.Cells.Item(3,4).Value = 0.04

$wsTotal.Cells.Item(4,1).Value = 2
$wsTotal.Cells.Item(4,2).Value = "2022-Q2"
$wsTotal.Cells.Item(4,3).Value = 1
$wsTotal.Cells.Item(4,4).Value = 0.03

$wsTotal.Cells.Item(5,1).Value = 3
$wsTotal.Cells.Item(5,2).Value = "2022-Q1"
$wsTotal.Cells.Item(5,3).Value = 1
$wsTotal.Cells.Item(5,4).Value = 0.03

# Restore the originally-active tab ("2022-Q1"); creating/copying sheets above
# shifts Excel's active-sheet selection around.
$wb.Worksheets.Item("2022-Q1").Activate()

Write-Output "done"
